$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 26, shifting the existing
# data (old rows 26-97) down to rows 28-99.
$ws.Rows("26:27").Insert()

# New row 26: "Especial" quality entry for the new reporting date.
$ws.Range("A26").Value = 1
$ws.Range("B26").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C26").Value = 'Arica y Parinacota'
$ws.Range("D26").Value = 44487
$ws.Range("E26").Value = 15
$ws.Range("F26").Value = 'Fruta'
$ws.Range("G26").Value = 100108
$ws.Range("H26").Value = 'Tropicales y subtropicales'
$ws.Range("I26").Value = 100108002
$ws.Range("J26").Value = 'Mango'
$ws.Range("K26").Value = 'Sin especificar'
$ws.Range("L26").Value = 'Especial'
$ws.Range("M26").Value = 450
$ws.Range("N26").Value = 4000
$ws.Range("O26").Value = 4500
$ws.Range("P26").Value = 4250
$ws.Range("Q26").Value = '$/bandeja 4 kilos'
$ws.Range("R26").Value = 'Perú'
$ws.Range("S26").Value = 1062
$ws.Range("T26").Value = 4

# New row 27: "Primera" quality entry for the same new reporting date.
$ws.Range("A27").Value = 1
$ws.Range("B27").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C27").Value = 'Arica y Parinacota'
$ws.Range("D27").Value = 44487
$ws.Range("E27").Value = 15
$ws.Range("F27").Value = 'Fruta'
$ws.Range("G27").Value = 100108
$ws.Range("H27").Value = 'Tropicales y subtropicales'
$ws.Range("I27").Value = 100108002
$ws.Range("J27").Value = 'Mango'
$ws.Range("K27").Value = 'Sin especificar'
$ws.Range("L27").Value = 'Primera'
$ws.Range("M27").Value = 450
$ws.Range("N27").Value = 4000
$ws.Range("O27").Value = 4500
$ws.Range("P27").Value = 4250
$ws.Range("Q27").Value = '$/bandeja 4 kilos'
$ws.Range("R27").Value = 'Perú'
$ws.Range("S27").Value = 1062
$ws.Range("T27").Value = 4
